$d = $word.ActiveDocument

function Set-ParagraphXml($paragraphIndex, $innerWml) {
    $p = $d.Paragraphs($paragraphIndex)
    $r = $p.Range
    $pkg = '<?xml version="1.0" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerWml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg)
}

# Paragraph 2: "{m:for v | null}" + moved bookmark, field markers removed.
$para2 = '<w:p>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>{m:</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">for v | </w:t></w:r>' +
    '<w:r><w:t>null}</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:t xml:space="preserve">    </w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>&lt;---</w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>The iteration variable types must be collections ([null]).</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml 2 $para2

# Paragraph 3: "{m:v.name}" field markers removed.
$nbsp = [char]0x00A0
$para3 = '<w:p>' +
    ('<w:r><w:t>name' + $nbsp + '=</w:t></w:r>') +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>{</w:t></w:r>' +
    '<w:r><w:t>m</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">:v.name}</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">    </w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>&lt;---</w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>No collection type for the iterator v</w:t></w:r>' +
    '<w:r><w:t>,</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml 3 $para3

# Paragraph 4: "{m:endfor}" field markers removed.
$para4 = '<w:p>' +
    '<w:r><w:t>{</w:t></w:r>' +
    '<w:r><w:t>m:</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">endfor}</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml 4 $para4

Write-Output "done"
